$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New GPS training-session rows are being appended below the existing
# data table (which currently ends at row 780). First clone the row-780
# formatting (date format in column B, centered style in column D) down
# across the full new block (781:798) so every new row matches the
# look of the existing rows, then fill in the actual values cell by cell.
$ws.Range("A780:V780").Copy()
$ws.Range("A781:V798").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 781
$ws.Cells.Item(781, 1).Value = "Entrainement"
$ws.Cells.Item(781, 2).Value = 45960
$ws.Cells.Item(781, 3).Value = "Global"
$ws.Cells.Item(781, 4).Value = "J-2"
$ws.Cells.Item(781, 5).Value = "Maé Clavel"
$ws.Cells.Item(781, 6).Value = "left back"
$ws.Cells.Item(781, 7).Value = "01:20:40"
$ws.Cells.Item(781, 8).Value = 4.25
$ws.Cells.Item(781, 9).Value = 0.36
$ws.Cells.Item(781, 10).Value = 3.89
$ws.Cells.Item(781, 11).Value = 0.2
$ws.Cells.Item(781, 12).Value = 0.1
$ws.Cells.Item(781, 13).Value = 0.06
$ws.Cells.Item(781, 14).Value = 0
$ws.Cells.Item(781, 15).Value = 4
$ws.Cells.Item(781, 16).Value = 2.95
$ws.Cells.Item(781, 17).Value = 28.05
$ws.Cells.Item(781, 18).Value = 4.42
$ws.Cells.Item(781, 19).Value = 9
$ws.Cells.Item(781, 20).Value = 3
$ws.Cells.Item(781, 21).Value = 5
$ws.Cells.Item(781, 22).Value = 0

# Row 782
$ws.Cells.Item(782, 1).Value = "Entrainement"
$ws.Cells.Item(782, 2).Value = 45960
$ws.Cells.Item(782, 3).Value = "Global"
$ws.Cells.Item(782, 4).Value = "J-2"
$ws.Cells.Item(782, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(782, 6).Value = "center midfield"
$ws.Cells.Item(782, 7).Value = "01:17:20"
$ws.Cells.Item(782, 8).Value = 4.68
$ws.Cells.Item(782, 9).Value = 0.51
$ws.Cells.Item(782, 10).Value = 4.1500000000000004
$ws.Cells.Item(782, 11).Value = 0.34
$ws.Cells.Item(782, 12).Value = 0.13
$ws.Cells.Item(782, 13).Value = 0.06
$ws.Cells.Item(782, 14).Value = 0
$ws.Cells.Item(782, 15).Value = 4
$ws.Cells.Item(782, 16).Value = 3.49
$ws.Cells.Item(782, 17).Value = 27.59
$ws.Cells.Item(782, 18).Value = 4.32
$ws.Cells.Item(782, 19).Value = 22
$ws.Cells.Item(782, 20).Value = 3
$ws.Cells.Item(782, 21).Value = 12
$ws.Cells.Item(782, 22).Value = 1

# Row 783
$ws.Cells.Item(783, 1).Value = "Entrainement"
$ws.Cells.Item(783, 2).Value = 45960
$ws.Cells.Item(783, 3).Value = "Global"
$ws.Cells.Item(783, 4).Value = "J-2"
$ws.Cells.Item(783, 5).Value = "Mattheo Haon"
$ws.Cells.Item(783, 6).Value = "right back"
$ws.Cells.Item(783, 7).Value = "01:17:19"
$ws.Cells.Item(783, 8).Value = 4.12
$ws.Cells.Item(783, 9).Value = 0.39
$ws.Cells.Item(783, 10).Value = 3.73
$ws.Cells.Item(783, 11).Value = 0.13
$ws.Cells.Item(783, 12).Value = 0.18
$ws.Cells.Item(783, 13).Value = 0.08
$ws.Cells.Item(783, 14).Value = 0
$ws.Cells.Item(783, 15).Value = 5
$ws.Cells.Item(783, 16).Value = 3.04
$ws.Cells.Item(783, 17).Value = 29.25
$ws.Cells.Item(783, 18).Value = 3.61
$ws.Cells.Item(783, 19).Value = 9
$ws.Cells.Item(783, 20).Value = 0
$ws.Cells.Item(783, 21).Value = 13
$ws.Cells.Item(783, 22).Value = 0

# Row 784
$ws.Cells.Item(784, 1).Value = "Entrainement"
$ws.Cells.Item(784, 2).Value = 45960
$ws.Cells.Item(784, 3).Value = "Global"
$ws.Cells.Item(784, 4).Value = "J-2"
$ws.Cells.Item(784, 5).Value = "Malik Boussaid"
$ws.Cells.Item(784, 6).Value = "right back"
$ws.Cells.Item(784, 7).Value = "01:17:51"
$ws.Cells.Item(784, 8).Value = 4.3499999999999996
$ws.Cells.Item(784, 9).Value = 0.31
$ws.Cells.Item(784, 10).Value = 4.04
$ws.Cells.Item(784, 11).Value = 0.12
$ws.Cells.Item(784, 12).Value = 0.12
$ws.Cells.Item(784, 13).Value = 0.06
$ws.Cells.Item(784, 14).Value = 0.01
$ws.Cells.Item(784, 15).Value = 6
$ws.Cells.Item(784, 16).Value = 2.82
$ws.Cells.Item(784, 17).Value = 31.11
$ws.Cells.Item(784, 18).Value = 4.59
$ws.Cells.Item(784, 19).Value = 16
$ws.Cells.Item(784, 20).Value = 7
$ws.Cells.Item(784, 21).Value = 14
$ws.Cells.Item(784, 22).Value = 13

# Row 785
$ws.Cells.Item(785, 1).Value = "Entrainement"
$ws.Cells.Item(785, 2).Value = 45960
$ws.Cells.Item(785, 3).Value = "Global"
$ws.Cells.Item(785, 4).Value = "J-2"
$ws.Cells.Item(785, 5).Value = "Sofiane Belle"
$ws.Cells.Item(785, 6).Value = "left forward"
$ws.Cells.Item(785, 7).Value = "01:14:55"
$ws.Cells.Item(785, 8).Value = 4.6100000000000003
$ws.Cells.Item(785, 9).Value = 0.37
$ws.Cells.Item(785, 10).Value = 4.24
$ws.Cells.Item(785, 11).Value = 0.22
$ws.Cells.Item(785, 12).Value = 0.12
$ws.Cells.Item(785, 13).Value = 0.04
$ws.Cells.Item(785, 14).Value = 0
$ws.Cells.Item(785, 15).Value = 4
$ws.Cells.Item(785, 16).Value = 3.52
$ws.Cells.Item(785, 17).Value = 26.65
$ws.Cells.Item(785, 18).Value = 4.26
$ws.Cells.Item(785, 19).Value = 23
$ws.Cells.Item(785, 20).Value = 5
$ws.Cells.Item(785, 21).Value = 16
$ws.Cells.Item(785, 22).Value = 1

# Row 786
$ws.Cells.Item(786, 1).Value = "Entrainement"
$ws.Cells.Item(786, 2).Value = 45960
$ws.Cells.Item(786, 3).Value = "Global"
$ws.Cells.Item(786, 4).Value = "J-2"
$ws.Cells.Item(786, 5).Value = "Yoan Zouma"
$ws.Cells.Item(786, 6).Value = "center back"
$ws.Cells.Item(786, 7).Value = "00:43:16"
$ws.Cells.Item(786, 8).Value = 2.33
$ws.Cells.Item(786, 9).Value = 0.14000000000000001
$ws.Cells.Item(786, 10).Value = 2.1800000000000002
$ws.Cells.Item(786, 11).Value = 0.04
$ws.Cells.Item(786, 12).Value = 0.09
$ws.Cells.Item(786, 13).Value = 0.02
$ws.Cells.Item(786, 14).Value = 0
$ws.Cells.Item(786, 15).Value = 4
$ws.Cells.Item(786, 16).Value = 2.96
$ws.Cells.Item(786, 17).Value = 26.08
$ws.Cells.Item(786, 18).Value = 3.25
$ws.Cells.Item(786, 19).Value = 2
$ws.Cells.Item(786, 20).Value = 0
$ws.Cells.Item(786, 21).Value = 0
$ws.Cells.Item(786, 22).Value = 1

# Row 787
$ws.Cells.Item(787, 1).Value = "Entrainement"
$ws.Cells.Item(787, 2).Value = 45960
$ws.Cells.Item(787, 3).Value = "Global"
$ws.Cells.Item(787, 4).Value = "J-2"
$ws.Cells.Item(787, 5).Value = "Karahali Souaré"
$ws.Cells.Item(787, 6).Value = "right forward"
$ws.Cells.Item(787, 7).Value = "01:00:22"
$ws.Cells.Item(787, 8).Value = 3.94
$ws.Cells.Item(787, 9).Value = 0.2
$ws.Cells.Item(787, 10).Value = 3.73
$ws.Cells.Item(787, 11).Value = 0.12
$ws.Cells.Item(787, 12).Value = 0.07
$ws.Cells.Item(787, 13).Value = 0.01
$ws.Cells.Item(787, 14).Value = 0
$ws.Cells.Item(787, 15).Value = 2
$ws.Cells.Item(787, 16).Value = 3.58
$ws.Cells.Item(787, 17).Value = 29
$ws.Cells.Item(787, 18).Value = 5.03
$ws.Cells.Item(787, 19).Value = 43
$ws.Cells.Item(787, 20).Value = 8
$ws.Cells.Item(787, 21).Value = 41
$ws.Cells.Item(787, 22).Value = 4

# Row 788
$ws.Cells.Item(788, 1).Value = "Entrainement"
$ws.Cells.Item(788, 2).Value = 45960
$ws.Cells.Item(788, 3).Value = "Global"
$ws.Cells.Item(788, 4).Value = "J-2"
$ws.Cells.Item(788, 5).Value = "Omar Benyounes"
$ws.Cells.Item(788, 6).Value = "center midfield"
$ws.Cells.Item(788, 7).Value = "01:19:03"
$ws.Cells.Item(788, 8).Value = 3.96
$ws.Cells.Item(788, 9).Value = 0.16
$ws.Cells.Item(788, 10).Value = 3.79
$ws.Cells.Item(788, 11).Value = 0.14000000000000001
$ws.Cells.Item(788, 12).Value = 0.03
$ws.Cells.Item(788, 13).Value = 0
$ws.Cells.Item(788, 14).Value = 0
$ws.Cells.Item(788, 15).Value = 0
$ws.Cells.Item(788, 16).Value = 2.78
$ws.Cells.Item(788, 17).Value = 23.67
$ws.Cells.Item(788, 18).Value = 4.1900000000000004
$ws.Cells.Item(788, 19).Value = 10
$ws.Cells.Item(788, 20).Value = 1
$ws.Cells.Item(788, 21).Value = 6
$ws.Cells.Item(788, 22).Value = 0

# Row 789
$ws.Cells.Item(789, 1).Value = "Entrainement"
$ws.Cells.Item(789, 2).Value = 45960
$ws.Cells.Item(789, 3).Value = "Global"
$ws.Cells.Item(789, 4).Value = "J-2"
$ws.Cells.Item(789, 5).Value = "Levy Ndoutoume"
$ws.Cells.Item(789, 6).Value = "left back"
$ws.Cells.Item(789, 7).Value = "01:14:07"
$ws.Cells.Item(789, 8).Value = 3.37
$ws.Cells.Item(789, 9).Value = 0.34
$ws.Cells.Item(789, 10).Value = 3.02
$ws.Cells.Item(789, 11).Value = 0.21
$ws.Cells.Item(789, 12).Value = 0.06
$ws.Cells.Item(789, 13).Value = 0.06
$ws.Cells.Item(789, 14).Value = 0.02
$ws.Cells.Item(789, 15).Value = 4
$ws.Cells.Item(789, 16).Value = 2.2000000000000002
$ws.Cells.Item(789, 17).Value = 31.21
$ws.Cells.Item(789, 18).Value = 4.96
$ws.Cells.Item(789, 19).Value = 13
$ws.Cells.Item(789, 20).Value = 5
$ws.Cells.Item(789, 21).Value = 13
$ws.Cells.Item(789, 22).Value = 8

# Row 790
$ws.Cells.Item(790, 1).Value = "Entrainement"
$ws.Cells.Item(790, 2).Value = 45960
$ws.Cells.Item(790, 3).Value = "Global"
$ws.Cells.Item(790, 4).Value = "J-2"
$ws.Cells.Item(790, 5).Value = "Naim Ighbane"
$ws.Cells.Item(790, 6).Value = "center back"
$ws.Cells.Item(790, 7).Value = "00:55:33"
$ws.Cells.Item(790, 8).Value = 4.78
$ws.Cells.Item(790, 9).Value = 1.41
$ws.Cells.Item(790, 10).Value = 3.36
$ws.Cells.Item(790, 11).Value = 0.8
$ws.Cells.Item(790, 12).Value = 0.59
$ws.Cells.Item(790, 13).Value = 0.04
$ws.Cells.Item(790, 14).Value = 0
$ws.Cells.Item(790, 15).Value = 4
$ws.Cells.Item(790, 16).Value = 4.6399999999999997
$ws.Cells.Item(790, 17).Value = 29.06
$ws.Cells.Item(790, 18).Value = 5.04
$ws.Cells.Item(790, 19).Value = 26
$ws.Cells.Item(790, 20).Value = 9
$ws.Cells.Item(790, 21).Value = 16
$ws.Cells.Item(790, 22).Value = 3

# Row 791
$ws.Cells.Item(791, 1).Value = "Entrainement"
$ws.Cells.Item(791, 2).Value = 45960
$ws.Cells.Item(791, 3).Value = "Global"
$ws.Cells.Item(791, 4).Value = "J-2"
$ws.Cells.Item(791, 5).Value = "Karim Belmahi"
$ws.Cells.Item(791, 6).Value = "left forward"
$ws.Cells.Item(791, 7).Value = "01:15:11"
$ws.Cells.Item(791, 8).Value = 4.38
$ws.Cells.Item(791, 9).Value = 0.3
$ws.Cells.Item(791, 10).Value = 4.08
$ws.Cells.Item(791, 11).Value = 0.16
$ws.Cells.Item(791, 12).Value = 0.12
$ws.Cells.Item(791, 13).Value = 0.02
$ws.Cells.Item(791, 14).Value = 0
$ws.Cells.Item(791, 15).Value = 6
$ws.Cells.Item(791, 16).Value = 3.01
$ws.Cells.Item(791, 17).Value = 26.38
$ws.Cells.Item(791, 18).Value = 5.19
$ws.Cells.Item(791, 19).Value = 18
$ws.Cells.Item(791, 20).Value = 11
$ws.Cells.Item(791, 21).Value = 19
$ws.Cells.Item(791, 22).Value = 5

# Row 792
$ws.Cells.Item(792, 1).Value = "Entrainement"
$ws.Cells.Item(792, 2).Value = 45961
$ws.Cells.Item(792, 3).Value = "Global"
$ws.Cells.Item(792, 4).Value = "J-1"
$ws.Cells.Item(792, 5).Value = "Karahali Souaré"
$ws.Cells.Item(792, 6).Value = "right forward"
$ws.Cells.Item(792, 7).Value = "00:34:44"
$ws.Cells.Item(792, 8).Value = 2.4500000000000002
$ws.Cells.Item(792, 9).Value = 0.22
$ws.Cells.Item(792, 10).Value = 2.2200000000000002
$ws.Cells.Item(792, 11).Value = 0.15
$ws.Cells.Item(792, 12).Value = 0.06
$ws.Cells.Item(792, 13).Value = 0.01
$ws.Cells.Item(792, 14).Value = 0
$ws.Cells.Item(792, 15).Value = 3
$ws.Cells.Item(792, 16).Value = 3.92
$ws.Cells.Item(792, 17).Value = 28.12
$ws.Cells.Item(792, 18).Value = 5.17
$ws.Cells.Item(792, 19).Value = 25
$ws.Cells.Item(792, 20).Value = 13
$ws.Cells.Item(792, 21).Value = 24
$ws.Cells.Item(792, 22).Value = 7

# Row 793
$ws.Cells.Item(793, 1).Value = "Entrainement"
$ws.Cells.Item(793, 2).Value = 45961
$ws.Cells.Item(793, 3).Value = "Global"
$ws.Cells.Item(793, 4).Value = "J-1"
$ws.Cells.Item(793, 5).Value = "Mattheo Haon"
$ws.Cells.Item(793, 6).Value = "right back"
$ws.Cells.Item(793, 7).Value = "00:36:22"
$ws.Cells.Item(793, 8).Value = 3.1
$ws.Cells.Item(793, 9).Value = 0.31
$ws.Cells.Item(793, 10).Value = 2.79
$ws.Cells.Item(793, 11).Value = 0.19
$ws.Cells.Item(793, 12).Value = 0.08
$ws.Cells.Item(793, 13).Value = 0.04
$ws.Cells.Item(793, 14).Value = 0
$ws.Cells.Item(793, 15).Value = 3
$ws.Cells.Item(793, 16).Value = 5.08
$ws.Cells.Item(793, 17).Value = 29.31
$ws.Cells.Item(793, 18).Value = 4.4400000000000004
$ws.Cells.Item(793, 19).Value = 24
$ws.Cells.Item(793, 20).Value = 4
$ws.Cells.Item(793, 21).Value = 20
$ws.Cells.Item(793, 22).Value = 1

# Row 794
$ws.Cells.Item(794, 1).Value = "Entrainement"
$ws.Cells.Item(794, 2).Value = 45961
$ws.Cells.Item(794, 3).Value = "Global"
$ws.Cells.Item(794, 4).Value = "J-1"
$ws.Cells.Item(794, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(794, 6).Value = "center midfield"
$ws.Cells.Item(794, 7).Value = "00:36:44"
$ws.Cells.Item(794, 8).Value = 3.35
$ws.Cells.Item(794, 9).Value = 0.13
$ws.Cells.Item(794, 10).Value = 3.21
$ws.Cells.Item(794, 11).Value = 0.13
$ws.Cells.Item(794, 12).Value = 0.01
$ws.Cells.Item(794, 13).Value = 0
$ws.Cells.Item(794, 14).Value = 0
$ws.Cells.Item(794, 15).Value = 0
$ws.Cells.Item(794, 16).Value = 5.43
$ws.Cells.Item(794, 17).Value = 21.35
$ws.Cells.Item(794, 18).Value = 4.0999999999999996
$ws.Cells.Item(794, 19).Value = 12
$ws.Cells.Item(794, 20).Value = 1
$ws.Cells.Item(794, 21).Value = 8
$ws.Cells.Item(794, 22).Value = 1

# Row 795
$ws.Cells.Item(795, 1).Value = "Entrainement"
$ws.Cells.Item(795, 2).Value = 45961
$ws.Cells.Item(795, 3).Value = "Global"
$ws.Cells.Item(795, 4).Value = "J-1"
$ws.Cells.Item(795, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(795, 6).Value = "left forward"
$ws.Cells.Item(795, 7).Value = "00:36:44"
$ws.Cells.Item(795, 8).Value = 3.33
$ws.Cells.Item(795, 9).Value = 0.31
$ws.Cells.Item(795, 10).Value = 3.01
$ws.Cells.Item(795, 11).Value = 0.21
$ws.Cells.Item(795, 12).Value = 0.07
$ws.Cells.Item(795, 13).Value = 0.03
$ws.Cells.Item(795, 14).Value = 0
$ws.Cells.Item(795, 15).Value = 5
$ws.Cells.Item(795, 16).Value = 5.08
$ws.Cells.Item(795, 17).Value = 30.14
$ws.Cells.Item(795, 18).Value = 4.6100000000000003
$ws.Cells.Item(795, 19).Value = 24
$ws.Cells.Item(795, 20).Value = 4
$ws.Cells.Item(795, 21).Value = 26
$ws.Cells.Item(795, 22).Value = 5

# Row 796
$ws.Cells.Item(796, 1).Value = "Entrainement"
$ws.Cells.Item(796, 2).Value = 45961
$ws.Cells.Item(796, 3).Value = "Global"
$ws.Cells.Item(796, 4).Value = "J-1"
$ws.Cells.Item(796, 5).Value = "Levy Ndoutoume"
$ws.Cells.Item(796, 6).Value = "left back"
$ws.Cells.Item(796, 7).Value = "00:36:07"
$ws.Cells.Item(796, 8).Value = 2.5
$ws.Cells.Item(796, 9).Value = 0.23
$ws.Cells.Item(796, 10).Value = 2.27
$ws.Cells.Item(796, 11).Value = 0.16
$ws.Cells.Item(796, 12).Value = 0.05
$ws.Cells.Item(796, 13).Value = 0.03
$ws.Cells.Item(796, 14).Value = 0
$ws.Cells.Item(796, 15).Value = 3
$ws.Cells.Item(796, 16).Value = 4.08
$ws.Cells.Item(796, 17).Value = 28.64
$ws.Cells.Item(796, 18).Value = 4.3499999999999996
$ws.Cells.Item(796, 19).Value = 23
$ws.Cells.Item(796, 20).Value = 4
$ws.Cells.Item(796, 21).Value = 10
$ws.Cells.Item(796, 22).Value = 3

# Row 797
$ws.Cells.Item(797, 1).Value = "Entrainement"
$ws.Cells.Item(797, 2).Value = 45961
$ws.Cells.Item(797, 3).Value = "Global"
$ws.Cells.Item(797, 4).Value = "J-1"
$ws.Cells.Item(797, 5).Value = "Maé Clavel"
$ws.Cells.Item(797, 6).Value = "left back"
$ws.Cells.Item(797, 7).Value = "00:35:37"
$ws.Cells.Item(797, 8).Value = 2.71
$ws.Cells.Item(797, 9).Value = 0.25
$ws.Cells.Item(797, 10).Value = 2.46
$ws.Cells.Item(797, 11).Value = 0.18
$ws.Cells.Item(797, 12).Value = 0.07
$ws.Cells.Item(797, 13).Value = 0.01
$ws.Cells.Item(797, 14).Value = 0
$ws.Cells.Item(797, 15).Value = 1
$ws.Cells.Item(797, 16).Value = 4.53
$ws.Cells.Item(797, 17).Value = 26.64
$ws.Cells.Item(797, 18).Value = 4.63
$ws.Cells.Item(797, 19).Value = 21
$ws.Cells.Item(797, 20).Value = 5
$ws.Cells.Item(797, 21).Value = 10
$ws.Cells.Item(797, 22).Value = 3

# Row 798
$ws.Cells.Item(798, 1).Value = "Entrainement"
$ws.Cells.Item(798, 2).Value = 45961
$ws.Cells.Item(798, 3).Value = "Global"
$ws.Cells.Item(798, 4).Value = "J-1"
$ws.Cells.Item(798, 5).Value = "Jeremie Laurent"
$ws.Cells.Item(798, 6).Value = "left forward"
$ws.Cells.Item(798, 7).Value = "00:34:59"
$ws.Cells.Item(798, 8).Value = 2.77
$ws.Cells.Item(798, 9).Value = 0.3
$ws.Cells.Item(798, 10).Value = 2.4700000000000002
$ws.Cells.Item(798, 11).Value = 0.2
$ws.Cells.Item(798, 12).Value = 0.1
$ws.Cells.Item(798, 13).Value = 0.01
$ws.Cells.Item(798, 14).Value = 0
$ws.Cells.Item(798, 15).Value = 2
$ws.Cells.Item(798, 16).Value = 4.7300000000000004
$ws.Cells.Item(798, 17).Value = 26.34
$ws.Cells.Item(798, 18).Value = 4.8600000000000003
$ws.Cells.Item(798, 19).Value = 28
$ws.Cells.Item(798, 20).Value = 9
$ws.Cells.Item(798, 21).Value = 19
$ws.Cells.Item(798, 22).Value = 5

# Match the workbooks on-screen selection state after scrolling down
# to review the newly appended rows.
$ws.Range("D802").Select()